$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user deleted several password entries (rows) from this sheet and
# added a new one at the bottom.
#
# Existing hyperlinks are anchored to fixed cell addresses and this
# engine does not auto-shift them when rows are deleted, so remove them
# all up front and re-create them afterwards at their new locations.
$ws.Hyperlinks.Delete()

# Delete the rows that were removed entirely. Work from the bottom up so
# earlier row numbers stay valid as later deletes happen.
# (original rows 2,3,4,5 and 15 are gone in the final sheet)
$ws.Rows(15).Delete()
$ws.Rows(5).Delete()
$ws.Rows(4).Delete()
$ws.Rows(3).Delete()
$ws.Rows(2).Delete()

# The "Desk ID" column (E) header/data is gone too - clear the leftover
# header cell so the used range shrinks back to A:D.
$ws.Range("E1").ClearContents()

# Add the new password entry as the last row (row 18 after the deletes
# above).
$ws.Range("A18").Value = "www.facebook.com"
$ws.Range("B18").Value = "https://www.facebook.com"
$ws.Range("C18").Value = 9367653559
$ws.Range("D18").Value = "mamapapaloveyou"

# Re-create the hyperlinks at their new (shifted) addresses.
$ws.Hyperlinks.Add($ws.Range("B11"), "https://accounts.google.com/")
$ws.Range("B11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:tumimbangyra@gmail.com")
$ws.Range("C13").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A14"), "http://www.facebook.com/")
$ws.Range("A14").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A15"), "http://www.facebook.com/")
$ws.Range("A15").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.facebook.com/")
$ws.Range("B15").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A16"), "http://www.facebook.com/")
$ws.Range("A16").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A17"), "http://www.streamlikers.cc/")
$ws.Range("A17").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B18"), "https://www.facebook.com/")
$ws.Range("B18").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A18"), "http://www.facebook.com/")
$ws.Range("A18").Style = "Hyperlink"

# Match the saved selection left by the user (row 2 header selected).
$null = $ws.Rows(2).Select()
